$d = $word.ActiveDocument

# 1. CHAT <...>: "CHAT_|Sender|..." -> "CHAT |Sender|..." (underscore becomes a space)
$d.Content.Find.Execute("CHAT_|Sender|", $true, $false, $false, $false, $false, $true, 1, $false, "CHAT |Sender|", 2) | Out-Null

# 2. Remove the stray "_GoBack" bookmark in the GAMST paragraph
$d.Bookmarks("_GoBack").Delete()

# 3. GETCA paragraph: add a trailing period
$d.Content.Find.Execute("eine Karte ziehen will", $true, $false, $false, $false, $false, $true, 1, $false, "eine Karte ziehen will.", 2) | Out-Null

# 4. THRCA paragraph: add a trailing period
$d.Content.Find.Execute("wegschmeissen (throw away) will", $true, $false, $false, $false, $false, $true, 1, $false, "wegschmeissen (throw away) will.", 2) | Out-Null

# 5. Insert a new "LOGOUT" paragraph right after the THRCA paragraph (paragraph 11)
$thrcaPara = $d.Paragraphs(11)
$thrcaPara.Range.InsertParagraphAfter()
$logoutPara = $d.Paragraphs(12)
$logoutPara.Range.InsertBefore("LOGOUT" + [char]9 + "Anfrage des Clients an den Server, aus dem Chat auszutreten.")

# 6. "(sockets)" -> "(Sockets)" and "Handlers. " -> "Handler. " in the Handler intro paragraph
$d.Content.Find.Execute("(sockets)", $true, $false, $false, $false, $false, $true, 1, $false, "(Sockets)", 2) | Out-Null
$d.Content.Find.Execute("Handlers. ", $true, $false, $false, $false, $false, $true, 1, $false, "Handler. ", 2) | Out-Null

# 7. "PingPongThread" -> "Thread"
$d.Content.Find.Execute("PingPongThread", $true, $false, $false, $false, $false, $true, 1, $false, "Thread", 2) | Out-Null
